# Insert a new data row at row 91 (pushing existing rows 91-185 down to 92-186)
# and populate it with the new record's values, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 91; this shifts rows 91-185 -> 92-186
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record
$ws.Range("A91").Value2 = 7
$ws.Range("B91").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value2 = "Ñuble"
$ws.Range("D91").Value2 = 44554
$ws.Range("E91").Value2 = 16
$ws.Range("F91").Value2 = 100112043
$ws.Range("G91").Value2 = "Pepino ensalada"
$ws.Range("H91").Value2 = "Sin especificar"
$ws.Range("I91").Value2 = "Primera"
$ws.Range("J91").Value2 = 200
$ws.Range("K91").Value2 = 7500
$ws.Range("L91").Value2 = 8000
$ws.Range("M91").Value2 = 7750
$ws.Range("N91").Value2 = "$/caja 80 unidades"
$ws.Range("O91").Value2 = "Región del Maule"
$ws.Range("P91").Value2 = 97
$ws.Range("Q91").Value2 = 80
$ws.Range("R91").Value2 = "Hortaliza"
